$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valve_2.0_600_1")

# --- Column C: factor 0.85 -> 0.9 (rows 2-12) ---
$ws.Range("C2:C12").Value = 0.9

# --- Column D: was a constant 0.6, now computed as C*C ---
# D2 is its own (non-shared) formula, D3:D12 becomes one shared-formula
# group, mirroring how column B is already laid out. Do this *before*
# touching column B below so the shared-formula group index ("si") lands
# on 1, matching the untouched B-column group's index of 0.
$ws.Range("D2").Formula = "=C2*C2"
$ws.Range("D3:D12").Formula = "=C3*C3"

# D12 kept the old style (index 25) while the rest of row 12 (B12/C12/G12)
# already use the thicker bottom-border style (index 26) -- copy that
# formatting over so D12 matches its row, without touching the formula
# we just wrote.
$ws.Range("C12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column B: multiplier 1.7 -> 4 ---
$ws.Range("B2:B11").Formula = "=G2*4"
$ws.Range("B12").Formula = "=G12*4"

# --- Selection moves from F20 to J16 ---
$ws.Activate()
$ws.Range("J16").Select()

# --- Rename sheet to flag the parabolic-plug variant ---
$ws.Name = "Valve_2.0_600_1(ParabolicPlug)"
